# Fixed weird string encoding issue
# Adds a new "col4" column (E) to Sheet1, filled with the constant value 5
# for every data row, mirroring the existing date/col1/col2/col3 columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("E1").Value = "col4"

# Fill E2:E23 with the value 5 (one row per existing data row)
$lastRow = 23
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = 5
}

# Match the author's final selection state (cell E23 selected)
$ws.Range("E23").Select()
